$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2' = '26.140.76'
    'E2' = '  +0.54%  '
    'D3' = '1.655.75'
    'E3' = '  +0.28%  '
    'D4' = '1.003'
    'E4' = '  +0.07%  '
    'D5' = '218.00'
    'E5' = '  +0.67%  '
    'D6' = '0.5297'
    'E6' = '  +1.88%  '
    'E7' = '  +0.04%  '
    'D8' = '0.2613'
    'E8' = '  -0.18%  '
    'D9' = '0.06332'
    'E9' = '  +1.18%  '
    'D10' = '20.44'
    'E10' = '  -0.60%  '
    'D11' = '0.07805'
    'E11' = '  +1.18%  '
    'D12' = '4.519'
    'E12' = '  +1.13%  '
    'D13' = '1.651.42'
    'E13' = '  -0.10%  '
    'D14' = '1.883.38'
    'E14' = '  +0.28%  '
    'E15' = '  +1.27%  '
    'D16' = '0.0₅8222'
    'E16' = '  +1.72%  '
    'D17' = '65.42'
    'E17' = '  +1.02%  '
    'D18' = '26.137.21'
    'E19' = '  +0.05%  '
    'D20' = '4.605'
    'E20' = '  +0.78%  '
    'D21' = '191.24'
    'E21' = '  +0.07%  '
    'E22' = '  +0.82%  '
    'E23' = '  +0.83%  '
    'E24' = '  +0.07%  '
    'D25' = '145.13'
    'E25' = '  +5.05%  '
    'D26' = '0.1230'
    'E26' = '  -0.15%  '
    'D27' = '7.219'
    'E27' = '  -0.30%  '
    'D28' = '15.99'
    'E28' = '  -0.36%  '
    'D29' = '1.462'
    'E29' = '  +4.37%  '
    'D30' = '0.05767'
    'E30' = '  -2.85%  '
    'D31' = '1.274'
    'E31' = '  +0.29%  '
    'D32' = '3.561'
    'E32' = '  +1.43%  '
    'E33' = '  +0.85%  '
    'D34' = '1.601'
    'E34' = '  +2.83%  '
    'D35' = '2.799'
    'E35' = '  +1.65%  '
    'D36' = '0.9514'
    'E36' = '  +0.32%  '
    'D37' = '2.415'
    'E37' = '  +0.07%  '
    'D38' = '0.5760'
    'E38' = '  +1.49%  '
    'E39' = '  +1.11%  '
    'D40' = '0.8541'
    'E40' = '  +1.00%  '
    'D41' = '5.790'
    'E41' = '  -1.62%  '
    'D42' = '104.72'
    'D43' = '1.044.31'
    'E43' = '  +4.41%  '
    'D45' = '1.797.46'
    'E45' = '  +0.20%  '
    'D46' = '56.89'
    'E46' = '  +0.51%  '
    'D47' = '1.002'
    'E47' = '  -0.36%  '
    'D48' = '0.4342'
    'E48' = '  +1.07%  '
    'D49' = '7.855'
    'E49' = '  -1.33%  '
    'D50' = '0.05151'
    'E50' = '  +0.02%  '
    'D51' = '1.444'
    'E51' = '  -1.86%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}
